$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update transition-probability matrix cells per updated game data (pulled March 7)
$ws.Cells.Item(2, 2).Value = 0.2154882154882155
$ws.Cells.Item(2, 3).Value = 0.5252525252525253
$ws.Cells.Item(2, 10).Value = 0.0101010101010101
$ws.Cells.Item(2, 16).Value = 0.1346801346801347
$ws.Cells.Item(2, 19).Value = 0.1144781144781145
$ws.Cells.Item(3, 2).Value = 0.006172839506172839
$ws.Cells.Item(3, 3).Value = 0.0308641975308642
$ws.Cells.Item(3, 10).Value = 0.02469135802469136
$ws.Cells.Item(3, 16).Value = 0.7098765432098766
$ws.Cells.Item(3, 19).Value = 0.228395061728395
$ws.Cells.Item(4, 16).Value = 0.7428571428571429
$ws.Cells.Item(4, 19).Value = 0.2571428571428571
$ws.Cells.Item(6, 2).Value = 0.06985294117647059
$ws.Cells.Item(6, 4).Value = 0.01102941176470588
$ws.Cells.Item(6, 6).Value = 0.07352941176470588
$ws.Cells.Item(6, 10).Value = 0.2647058823529412
$ws.Cells.Item(6, 15).Value = 0.01838235294117647
$ws.Cells.Item(6, 17).Value = 0.125
$ws.Cells.Item(6, 18).Value = 0.05147058823529412
$ws.Cells.Item(6, 19).Value = 0.3860294117647059
$ws.Cells.Item(7, 2).Value = 0.09417040358744394
$ws.Cells.Item(7, 4).Value = 0.02242152466367713
$ws.Cells.Item(7, 6).Value = 0.07623318385650224
$ws.Cells.Item(7, 10).Value = 0.1524663677130045
$ws.Cells.Item(7, 15).Value = 0.01345291479820628
$ws.Cells.Item(7, 17).Value = 0.1210762331838565
$ws.Cells.Item(7, 18).Value = 0.05829596412556054
$ws.Cells.Item(7, 19).Value = 0.4618834080717489
$ws.Cells.Item(8, 2).Value = 0.1154684095860566
$ws.Cells.Item(8, 4).Value = 0.01742919389978214
$ws.Cells.Item(8, 5).Value = 0.002178649237472767
$ws.Cells.Item(8, 6).Value = 0.08061002178649238
$ws.Cells.Item(8, 10).Value = 0.10239651416122
$ws.Cells.Item(8, 15).Value = 0.0196078431372549
$ws.Cells.Item(8, 17).Value = 0.1699346405228758
$ws.Cells.Item(8, 18).Value = 0.08932461873638345
$ws.Cells.Item(8, 19).Value = 0.4030501089324618
$ws.Cells.Item(9, 2).Value = 0.07526881720430108
$ws.Cells.Item(9, 4).Value = 0.01075268817204301
$ws.Cells.Item(9, 5).Value = 0.005376344086021506
$ws.Cells.Item(9, 6).Value = 0.08602150537634409
$ws.Cells.Item(9, 10).Value = 0.09677419354838709
$ws.Cells.Item(9, 15).Value = 0.01612903225806452
$ws.Cells.Item(9, 17).Value = 0.1881720430107527
$ws.Cells.Item(9, 18).Value = 0.05913978494623656
$ws.Cells.Item(9, 19).Value = 0.4623655913978494
$ws.Cells.Item(10, 2).Value = 0.0992616899097621
$ws.Cells.Item(10, 4).Value = 0.014766201804758
$ws.Cells.Item(10, 6).Value = 0.08285479901558655
$ws.Cells.Item(10, 10).Value = 0.118129614438064
$ws.Cells.Item(10, 15).Value = 0.01312551271534044
$ws.Cells.Item(10, 17).Value = 0.1862182116488925
$ws.Cells.Item(10, 18).Value = 0.09515996718621821
$ws.Cells.Item(10, 19).Value = 0.3904840032813782
$ws.Cells.Item(11, 7).Value = 0.1373134328358209
$ws.Cells.Item(11, 10).Value = 0.09552238805970149
$ws.Cells.Item(11, 11).Value = 0.1850746268656716
$ws.Cells.Item(11, 12).Value = 0.5611940298507463
$ws.Cells.Item(11, 19).Value = 0.0208955223880597
$ws.Cells.Item(12, 7).Value = 0.7731958762886598
$ws.Cells.Item(12, 10).Value = 0.1391752577319588
$ws.Cells.Item(12, 11).Value = 0.0154639175257732
$ws.Cells.Item(12, 12).Value = 0.03092783505154639
$ws.Cells.Item(12, 19).Value = 0.04123711340206185
$ws.Cells.Item(13, 7).Value = 0.6666666666666666
$ws.Cells.Item(13, 10).Value = 0.1875
$ws.Cells.Item(13, 19).Value = 0.1458333333333333
$ws.Cells.Item(15, 6).Value = 0.01204819277108434
$ws.Cells.Item(15, 8).Value = 0.1646586345381526
$ws.Cells.Item(15, 9).Value = 0.06024096385542169
$ws.Cells.Item(15, 10).Value = 0.3895582329317269
$ws.Cells.Item(15, 11).Value = 0.09236947791164658
$ws.Cells.Item(15, 13).Value = 0.004016064257028112
$ws.Cells.Item(15, 15).Value = 0.07630522088353414
$ws.Cells.Item(15, 19).Value = 0.2008032128514056
$ws.Cells.Item(16, 6).Value = 0.04624277456647399
$ws.Cells.Item(16, 8).Value = 0.1849710982658959
$ws.Cells.Item(16, 9).Value = 0.06936416184971098
$ws.Cells.Item(16, 10).Value = 0.3410404624277457
$ws.Cells.Item(16, 11).Value = 0.1213872832369942
$ws.Cells.Item(16, 13).Value = 0.01734104046242774
$ws.Cells.Item(16, 15).Value = 0.05780346820809248
$ws.Cells.Item(16, 19).Value = 0.161849710982659
$ws.Cells.Item(17, 6).Value = 0.02238805970149254
$ws.Cells.Item(17, 8).Value = 0.1691542288557214
$ws.Cells.Item(17, 9).Value = 0.1044776119402985
$ws.Cells.Item(17, 10).Value = 0.3706467661691542
$ws.Cells.Item(17, 11).Value = 0.1044776119402985
$ws.Cells.Item(17, 13).Value = 0.02985074626865672
$ws.Cells.Item(17, 15).Value = 0.07213930348258707
$ws.Cells.Item(17, 19).Value = 0.1268656716417911
$ws.Cells.Item(18, 6).Value = 0.04060913705583756
$ws.Cells.Item(18, 8).Value = 0.1573604060913706
$ws.Cells.Item(18, 9).Value = 0.09137055837563451
$ws.Cells.Item(18, 10).Value = 0.3959390862944163
$ws.Cells.Item(18, 11).Value = 0.08121827411167512
$ws.Cells.Item(18, 13).Value = 0.02030456852791878
$ws.Cells.Item(18, 15).Value = 0.1065989847715736
$ws.Cells.Item(18, 19).Value = 0.1065989847715736
$ws.Cells.Item(19, 6).Value = 0.01864280387770321
$ws.Cells.Item(19, 8).Value = 0.2155108128262491
$ws.Cells.Item(19, 9).Value = 0.07606263982102908
$ws.Cells.Item(19, 10).Value = 0.3482475764354959
$ws.Cells.Item(19, 11).Value = 0.1230425055928412
$ws.Cells.Item(19, 13).Value = 0.02162565249813572
$ws.Cells.Item(19, 14).Value = 0.0007457121551081282
$ws.Cells.Item(19, 15).Value = 0.07979120059656973
$ws.Cells.Item(19, 19).Value = 0.116331096196868
